$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RMSE result value in A1 to reflect the recalculated
# mechanism coordinates (BCEF -> BCFG and CDGI -> CDEI) after porting
# the PMKS positionSolver to MATLAB.
$ws.Range("A1").Value = 77.110616148552836
